$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01281184547501269
$ws.Range("C2").Value = 0.2052201710742377

$ws.Range("B3").Value = 0.06977971597469282
$ws.Range("C3").Value = 0.1934796635155446

$ws.Range("B4").Value = 0.6876801604990849
$ws.Range("C4").Value = 0.158684699782748

$ws.Range("B5").Value = 0.9608464093066923
$ws.Range("C5").Value = 0.4135820600011968

$ws.Range("B6").Value = 0.9098575711074481
$ws.Range("C6").Value = 0.5373737244053245

$ws.Range("B7").Value = 0.7478133521552075
$ws.Range("C7").Value = 0.08705150454388995

$ws.Range("B8").Value = 0.006505982875823975
$ws.Range("C8").Value = 0.212273006439209
